# Append a new data row (row 67) to each of the four sheets, mirroring the
# schema of the existing rows (time, length, ID, actual length, checksum,
# and their decimal counterparts).

$wb = $excel.ActiveWorkbook

$rows = @{
    "MID_LFT_#1" = @{
        A = 45853.45988425926
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x50"
        E = "0x07"
        F = 400
        G = "5.68631262647113e+23"
        H = 336
        I = 7
    }
    "MID_LFT_#2" = @{
        A = 45853.45988425926
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x50"
        E = "0x19"
        F = 380
        G = "5.68432987514711e+23"
        H = 336
        I = 25
    }
    "MID_PLT_#1" = @{
        A = 45853.45988425926
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x66"
        E = "0x15"
        F = 110
        G = "5.68631262647113e+23"
        H = 102
        I = 15
    }
    "MID_PLT_#2" = @{
        A = 45853.45988425926
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x7A"
        E = "0x9"
        F = 130
        G = "5.68631262647113e+23"
        H = 122
        I = 9
    }
}

foreach ($ws in $wb.Worksheets) {
    $data = $rows[$ws.Name]
    if ($data -eq $null) { continue }

    $newRow = 67

    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = [double]$data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
